$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new blank rows above the existing row 3, pushing its data down to row 8
$ws.Range("A3:A7").EntireRow.Insert()
$ws.Range("A3:K7").ClearFormats()

# Row 3
$ws.Range("A3").Value = 1
$ws.Range("A3").NumberFormat = "0"
$ws.Range("B3").Value = "13-12-2025"
$ws.Range("C3").Value = "pranav"
$ws.Range("D3").Value = 9746271355
$ws.Range("D3").NumberFormat = "0"
$ws.Range("E3").Value = "21-12-2025"
$ws.Range("F3").Value = "ARJUN P"
$ws.Range("G3").Value = "Loss"
$ws.Range("H3").Value = "ENQUIRY"
$ws.Range("I3").Value = "Enquiry for Relative/Friend"
$ws.Range("J3").Value = "-"
$ws.Range("K3").Value = "just checking"

# Row 4
$ws.Range("A4").Value = 2
$ws.Range("A4").NumberFormat = "0"
$ws.Range("B4").Value = "13-12-2025"
$ws.Range("C4").Value = "SAYOOJ"
$ws.Range("D4").Value = 9539251325
$ws.Range("D4").NumberFormat = "0"
$ws.Range("E4").Value = "21-12-2025"
$ws.Range("F4").Value = "ARJUN P"
$ws.Range("G4").Value = "Loss"
$ws.Range("H4").Value = "ENQUIRY"
$ws.Range("I4").Value = "Enquiry for Relative/Friend"
$ws.Range("J4").Value = "-"
$ws.Range("K4").Value = "just checking"

# Row 5
$ws.Range("A5").Value = 3
$ws.Range("A5").NumberFormat = "0"
$ws.Range("B5").Value = "14-12-2025"
$ws.Range("C5").Value = "RAHANAS"
$ws.Range("D5").Value = 8848991304
$ws.Range("D5").NumberFormat = "0"
$ws.Range("E5").Value = "17-12-2025"
$ws.Range("F5").Value = "ARJUN P"
$ws.Range("G5").Value = "Loss"
$ws.Range("H5").Value = "CUSTOMER INTERNAL ISSUES"
$ws.Range("I5").Value = "FAMILY DISAPPROVEL"
$ws.Range("J5").Value = "-"
$ws.Range("K5").Value = "Tommorow coming"

# Row 6
$ws.Range("A6").Value = 4
$ws.Range("A6").NumberFormat = "0"
$ws.Range("B6").Value = "16-12-2025"
$ws.Range("C6").Value = "munaver"
$ws.Range("D6").Value = 8590020444
$ws.Range("D6").NumberFormat = "0"
$ws.Range("E6").Value = "25-12-2025"
$ws.Range("F6").Value = "ARJUN P"
$ws.Range("G6").Value = "Loss"
$ws.Range("H6").Value = "SIZE NOT SUITABLE"
$ws.Range("I6").Value = "SIZE TOO SMALL"
$ws.Range("J6").Value = "-"
$ws.Range("K6").Value = "size problem contact in two days"

# Row 7
$ws.Range("A7").Value = 5
$ws.Range("A7").NumberFormat = "0"
$ws.Range("B7").Value = "16-12-2025"
$ws.Range("C7").Value = "ATHUL"
$ws.Range("D7").Value = 9061301868
$ws.Range("D7").NumberFormat = "0"
$ws.Range("E7").Value = "28-12-2025"
$ws.Range("F7").Value = "ARJUN P"
$ws.Range("G7").Value = "Loss"
$ws.Range("H7").Value = "ENQUIRY"
$ws.Range("I7").Value = "ENQUIRY WITHOUT TRIAL"
$ws.Range("J7").Value = "-"
$ws.Range("K7").Value = "just checking"

# The original entry (Praveen) shifted from row 3 down to row 8 automatically via Insert();
# renumber its "#" value to keep the sequence 1..6
$ws.Range("A8").Value = 6
$ws.Range("A8").NumberFormat = "0"

# Widen the "Category" (H) and "Remarks" (K) columns to fit the new, longer text
$ws.Columns.Item(8).ColumnWidth = 31.5
$ws.Columns.Item(11).ColumnWidth = 42.333333333333336
